$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion note text (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")
$oldText = $cellA1.Value2
$newText = $oldText -replace [regex]::Escape("1000 Bs = 2.6 = 9623.38 pesos"), "1000 Bs = 2.68 = 9929.96 pesos"
$newText = $newText -replace [regex]::Escape("9623.38 pesos = 2.58 = 927.51 Bs"), "9929.96 pesos = 2.67 = 944.92 Bs"
$cellA1.Value2 = $newText

# --- Sheet "tasas": update the rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 372.811
$wsTasas.Range("O10").Value = 3702
$wsTasas.Range("N12").Value = 3720
$wsTasas.Range("O12").Value = 353.989
